$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 7 into row 8, then fill in new row values
$ws.Range("A7:E7").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)

# New row 8 data
$ws.Range("A8").Value = 43917.98333636574
$ws.Range("B8").Value = 123456.0
$ws.Range("C8").Value = "Qurator: Alternativ 2"
$ws.Range("D8").Value = "PQE: Vakant"
$ws.Range("E8").Value = "Novischförman: Alternativ 1, Novischförman: Alternativ 2, Novischförman: Blankt"
